$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly refresh of daily price-listing rows (Pepino dulce, Agricola del Norte).
# Each element: Row, D(Fecha-serial), H(Variedad), I(Calidad), J(Volumen),
# K(Precio minimo), L(Precio maximo), M(Precio promedio ponderado),
# N(Unidad de comercializacion), O(Origen), P(Precio $/Kg), Q(Kg o Unidades)
$rows = @(
    @(2, 44533, 'Cultivar XV región', 'Primera', 100, 6000, 7000, 6500, '$/caja 10 kilos', 'Región de Arica y Parinacota', 650, 10),
    @(3, 44533, 'Cultivar XV región', 'Segunda', 120, 4000, 5000, 4500, '$/caja 10 kilos', 'Región de Arica y Parinacota', 450, 10),
    @(4, 44211, 'Cultivar XV región', 'Segunda', 140, 4500, 5000, 4750, '$/caja 10 kilos', 'Región de Arica y Parinacota', 475, 10),
    @(5, 44391, 'Cultivar IV Región', 'Segunda', 100, 15000, 16000, 15500, '$/bandeja 18 kilos', 'Provincia de Limarí', 861, 18),
    @(6, 44554, 'Cultivar XV región', 'Primera', 200, 5000, 6000, 5500, '$/caja 10 kilos', 'Región de Arica y Parinacota', 550, 10),
    @(7, 44769, 'Cultivar IV Región', 'Primera', 140, 17000, 18000, 17500, '$/bandeja 18 kilos', 'Provincia de Limarí', 972, 18),
    @(8, 44377, 'Cultivar IV Región', 'Primera', 100, 17000, 18000, 17600, '$/bandeja 18 kilos', 'Provincia de Limarí', 978, 18),
    @(9, 44755, 'Cultivar IV Región', 'Primera', 160, 17000, 18000, 17500, '$/bandeja 18 kilos', 'Provincia de Limarí', 972, 18),
    @(10, 44748, 'Cultivar IV Región', 'Primera', 250, 17000, 18000, 17500, '$/bandeja 18 kilos', 'Provincia de Limarí', 972, 18),
    @(11, 44742, 'Cultivar IV Región', 'Segunda', 250, 15000, 16000, 15500, '$/bandeja 18 kilos', 'Provincia de Limarí', 861, 18),
    @(12, 45021, 'Cultivar IV Región', 'Primera', 270, 17000, 18000, 17500, '$/bandeja 18 kilos', 'Provincia de Limarí', 972, 18),
    @(13, 44783, 'Cultivar IV Región', 'Primera', 150, 17000, 18000, 17500, '$/bandeja 18 kilos', 'Provincia de Limarí', 972, 18),
    @(14, 44757, 'Cultivar XV región', 'Primera', 150, 6000, 6500, 6250, '$/caja 10 kilos', 'Región de Arica y Parinacota', 625, 10),
    @(15, 44433, 'Cultivar IV Región', 'Segunda', 100, 17000, 18000, 17500, '$/bandeja 18 kilos', 'Provincia de Limarí', 972, 18),
    @(16, 44433, 'Cultivar IV Región', 'Tercera', 120, 14000, 15000, 14500, '$/bandeja 18 kilos', 'Provincia de Limarí', 806, 18),
    @(17, 45035, 'Cultivar IV Región', 'Primera', 250, 19000, 20000, 19500, '$/bandeja 18 kilos', 'Provincia de Limarí', 1083, 18),
    @(18, 45042, 'Cultivar IV Región', 'Segunda', 220, 17000, 18000, 17545, '$/bandeja 18 kilos', 'Provincia de Limarí', 975, 18),
    @(19, 44405, 'Cultivar IV Región', 'Segunda', 140, 17000, 18000, 17500, '$/bandeja 18 kilos', 'Provincia de Limarí', 972, 18),
    @(20, 44412, 'Cultivar IV Región', 'Primera', 150, 17000, 18000, 17500, '$/bandeja 18 kilos', 'Provincia de Limarí', 972, 18),
    @(21, 44776, 'Cultivar IV Región', 'Primera', 200, 17000, 18000, 17500, '$/bandeja 18 kilos', 'Provincia de Limarí', 972, 18),
    @(22, 44771, 'Cultivar XV región', 'Primera', 140, 8000, 9000, 8500, '$/caja 10 kilos', 'Región de Arica y Parinacota', 850, 10),
    @(23, 44363, 'Cultivar IV Región', 'Primera', 140, 14000, 15000, 14500, '$/bandeja 18 kilos', 'Provincia de Limarí', 806, 18),
    @(24, 44762, 'Cultivar IV Región', 'Primera', 160, 15000, 16000, 15500, '$/bandeja 18 kilos', 'Provincia de Limarí', 861, 18),
    @(25, 44221, 'Cultivar XV región', 'Primera', 140, 5000, 6000, 5500, '$/caja 10 kilos', 'Región de Arica y Parinacota', 550, 10),
    @(26, 45114, 'Cultivar XV región', 'Primera', 160, 5000, 6000, 5500, '$/caja 10 kilos', 'Región de Arica y Parinacota', 550, 10),
    @(27, 44398, 'Cultivar IV Región', 'Primera', 100, 17000, 18000, 17500, '$/bandeja 18 kilos', 'Provincia de Limarí', 972, 18),
    @(28, 44398, 'Cultivar IV Región', 'Segunda', 100, 15000, 16000, 15500, '$/bandeja 18 kilos', 'Provincia de Limarí', 861, 18),
    @(29, 44435, 'Cultivar IV Región', 'Segunda', 100, 17000, 18000, 17500, '$/bandeja 18 kilos', 'Provincia de Limarí', 972, 18),
    @(30, 44435, 'Cultivar IV Región', 'Tercera', 120, 14000, 15000, 14500, '$/bandeja 18 kilos', 'Provincia de Limarí', 806, 18),
    @(31, 44454, 'Cultivar IV Región', 'Primera', 160, 19000, 20000, 19500, '$/bandeja 18 kilos', 'Provincia de Limarí', 1083, 18),
    @(32, 45043, 'Cultivar IV Región', 'Segunda', 170, 18000, 20000, 19059, '$/bandeja 18 kilos', 'Provincia de Limarí', 1059, 18),
    @(33, 44526, 'Cultivar XV región', 'Primera', 100, 5000, 5500, 5250, '$/caja 10 kilos', 'Región de Arica y Parinacota', 525, 10),
    @(34, 44526, 'Cultivar XV región', 'Segunda', 100, 4000, 4500, 4250, '$/caja 10 kilos', 'Región de Arica y Parinacota', 425, 10),
    @(35, 44526, 'Cultivar XV región', 'Tercera', 120, 3000, 3500, 3250, '$/caja 10 kilos', 'Región de Arica y Parinacota', 325, 10)
)

foreach ($r in $rows) {
    $row = $r[0]
    $ws.Cells.Item($row, 4).Value = $r[1]
    $ws.Cells.Item($row, 8).Value = $r[2]
    $ws.Cells.Item($row, 9).Value = $r[3]
    $ws.Cells.Item($row, 10).Value = $r[4]
    $ws.Cells.Item($row, 11).Value = $r[5]
    $ws.Cells.Item($row, 12).Value = $r[6]
    $ws.Cells.Item($row, 13).Value = $r[7]
    $ws.Cells.Item($row, 14).Value = $r[8]
    $ws.Cells.Item($row, 15).Value = $r[9]
    $ws.Cells.Item($row, 16).Value = $r[10]
    $ws.Cells.Item($row, 17).Value = $r[11]
}
